$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "휴림로봇"
$ws.Range("A6").Value = "엑스페릭스"
$ws.Range("A7").Value = "화신"
$ws.Range("A8").Value = "HLB테라퓨틱스"
$ws.Range("A9").Value = "에토니모리"

$ws.Range("A10").Select()
